# ---------------------------------------------------------------------------
# Adds a "Measures" sheet (DAX measure catalogue) and restructures "Data
# Insights" into a two-column (label / text) report layout.
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Data Insights")
$ws1.Activate()

# --- create "Measures" by duplicating "Data Insights", then wipe it out -----
# (duplicating keeps the same modern sheetFormatPr/namespace shape Excel uses;
# a blank Worksheets.Add() sheet would come out looking like a legacy sheet.)
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item($ws1.Index + 1)
$ws2.Name = "Measures"
$ws2.Cells.ClearContents()

# --- Sheet 1: "Data Insights" -----------------------------------------------
$ws1.Cells.ClearContents()

$ws1.Range("A1").Value = "Report Name"
$ws1.Range("B1").Value = "Data Insights"
$ws1.Range("A2").Value = "Sales Analysis"
$ws1.Range("B2").Value = "People are more interested in buying technology products"
$ws1.Range("B3").Value = "Phones are most bought sub category "
$ws1.Range("B4").Value = "Art supplies have one of the least sales. We need to motivate people to pursue Art."
$ws1.Range("B5").Value = "Country having most sales is United States"
$ws1.Range("B6").Value = "State having most sales is England"
$ws1.Range("B7").Value = "City having most sales is New york"
$ws1.Range("B8").Value = "Consumer products have the most sales"
$ws1.Range("B9").Value = "Most sales are recorded in Western Europe region followed by Central America"
$ws1.Range("A11").Value = "Profit Analysis"
$ws1.Range("B11").Value = "Technology products are most profitable category"
$ws1.Range("B12").Value = "Copiers are most profitable sub category"
$ws1.Range("B13").Value = "Australia and France have most sales than China but China is more profitable than them. Need to investigate why ?"
$ws1.Range("B14").Value = "India and UK have less sales comparatively but are in top 4 in terms of profits. Need to investigate why ?"
$ws1.Range("B15").Value = "Manila city has good sales but profits are less, need to investigate ?"
$ws1.Range("B16").Value = "Europe market is highly profitable and next is Asia pacific"
$ws1.Range("B17").Value = "Western Africa, Central Asia and Western Asia have negative profits"
$ws1.Range("B18").Value = "Global store sales are most profitable in United States"
$ws1.Range("B19").Value = "State where sales are most profitable is in England, next is California"
$ws1.Range("B20").Value = "Sales are most profitable in New york city"

# Header row (A1:B1) reuses the bold/fill/border "header" style already on A1.
$ws1.Range("A1").Copy()
$ws1.Range("B1").PasteSpecial(-4122)
$ws1.Application.CutCopyMode = $false

$ws1.Columns.Item(1).ColumnWidth = 11.666666666666666
$ws1.Columns.Item(2).ColumnWidth = 95

[void]$ws1.Range("B16").Select()

# --- Sheet 2: "Measures" ------------------------------------------------------
$ws2.Range("A1").Value = "Name"
$ws2.Range("B1").Value = "DAX Calculation"
$ws2.Range("A2").Value = "% of Sales by Region"
$ws2.Range("B2").Value = "DIVIDE([TotalSales], CALCULATE(SUM(Orders[Sales]), ALL(Orders[Region]))) * 100"
$ws2.Range("A4").Value = "% of Profits By Region"
$ws2.Range("B4").Value = "DIVIDE([TotalProfit], CALCULATE(SUM(Orders[Profit]), ALL(Orders[Region]))) * 100"
$ws2.Range("A6").Value = "Total Customers"
$ws2.Range("B6").Value = "DISTINCTCOUNT(Orders[Customer ID])"
$ws2.Range("A8").Value = "Total Orders"
$ws2.Range("B8").Value = "DISTINCTCOUNT(Orders[Order ID])"
$ws2.Range("A10").Value = "Total Sales"
$ws2.Range("B10").Value = "SUM(Orders[Sales])"
$ws2.Range("A12").Value = "Total Profit"
$ws2.Range("B12").Value = "SUM(Orders[Profit])"
$ws2.Range("A14").Value = "Total Shipping Cost"
$ws2.Range("B14").Value = "SUM(Orders[Shipping Cost])"

$ws1.Range("A1").Copy()
$ws2.Range("A1:B1").PasteSpecial(-4122)
$ws2.Application.CutCopyMode = $false

$ws2.Columns.Item(1).ColumnWidth = 18.333333333333332
$ws2.Columns.Item(2).ColumnWidth = 85.16666666666667

[void]$ws2.Range("C11").Select()

# Restore "Data Insights" as the active/selected tab (matches original file).
$ws1.Activate()
Write-Host "Workbook restructured: Data Insights (A:B) + Measures sheet added"
